$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
    "B1" = "product_name"
    "C1" = "productid"
    "D1" = "quality"
    "E1" = "total_quantity_ordered"
    "F1" = "total_revenue"
    "G1" = "thca_percentage"
    "H1" = "total_cbd"
    "I1" = "cbga"
    "J1" = "total_cbg"
    "K1" = "delta_nine_thc"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}
